$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B14").Value = "I need some other alternatives. I could use some other circuit characteristics. Rise time? Power? Predicted operating time still seems like a good one though."
